$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview": update status/date for the b.md row (row 3) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 02:37:01"

# ---- Sheet "zh-cn": update handoff info for the b.md row (row 3) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-23 02:36:56"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9acfc3f747f1e7f6324698189d0c983ed77e3c34/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8154847c6f0b693f50a0cbe9bc6b50f0d31ac54d/e2e/b.md."
# Widen the "Error Detail" column (P) to match the other long-text columns
$wsZhCn.Columns.Item(16).ColumnWidth = $wsZhCn.Columns.Item(7).ColumnWidth

# ---- Sheet "de-de": update handoff info for the b.md row (row 3) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-23 02:37:01"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9acfc3f747f1e7f6324698189d0c983ed77e3c34/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8154847c6f0b693f50a0cbe9bc6b50f0d31ac54d/e2e/b.md."
# Widen the "Error Detail" column (P) to match the other long-text columns
$wsDeDe.Columns.Item(16).ColumnWidth = $wsDeDe.Columns.Item(7).ColumnWidth
